$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old contents first (sheet had data through C23)
$ws.Range("A1:C23").ClearContents()

# New, reduced set of rows/values
$ws.Range("A1").Value = "Teacher"
$ws.Range("A2").Value = "a"
$ws.Range("A3").Value = "a"
$ws.Range("A5").Value = "Student"
$ws.Range("A8").Value = "Students Average:"
$ws.Range("B8").Value = 0
$ws.Range("A9").Value = "Median:"
$ws.Range("B9").Value = 0
$ws.Range("A10").Value = "Excellent Students:"
